$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to be treated as text so numeric-looking strings
# (e.g. "109.93") are stored as exact text, matching the original
# inlineStr cell type, instead of being auto-converted to floating
# point numbers with rounding artifacts.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.355.01"
$ws.Range("E2").Value = "  -3.04%  "

$ws.Range("D3").Value = "2.221.48"
$ws.Range("E3").Value = "  -2.32%  "

$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").Value = "109.93"
$ws.Range("E5").Value = "  -7.61%  "

$ws.Range("D6").Value = "287.03"
$ws.Range("E6").Value = "  +7.36%  "

$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -3.21%  "

$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("D9").Value = "0.597"
$ws.Range("E9").Value = "  -4.10%  "

$ws.Range("D10").Value = "43.31"
$ws.Range("E10").Value = "  -8.75%  "

$ws.Range("D11").Value = "0.0909"
$ws.Range("E11").Value = "  -3.87%  "

$ws.Range("D12").Value = "54.19"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").Value = "8.60"
$ws.Range("E13").Value = "  -9.10%  "

$ws.Range("D14").Value = "1.01"
$ws.Range("E14").Value = "  +12.56%  "

$ws.Range("E15").Value = "  -3.01%  "

$ws.Range("D16").Value = "14.82"
$ws.Range("E16").Value = "  -5.79%  "

$ws.Range("D17").Value = "2.551.58"
$ws.Range("E17").Value = "  -2.39%  "

$ws.Range("D18").Value = "2.233.97"

$ws.Range("D19").Value = "42.324.21"
$ws.Range("E19").Value = "  -2.92%  "

$ws.Range("D20").Value = "7.14"
$ws.Range("E20").Value = "  +3.06%  "

$ws.Range("D21").Value = "0.0000104"
$ws.Range("E21").Value = "  -4.88%  "

$ws.Range("D22").Value = "72.86"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").Value = "3.34"
$ws.Range("E23").Value = "  +14.14%  "

$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -0.99%  "

$ws.Range("D25").Value = "228.95"
$ws.Range("E25").Value = "  -2.39%  "

$ws.Range("D26").Value = "8.94"
$ws.Range("E26").Value = "  -6.56%  "

$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("D28").Value = "11.40"
$ws.Range("E28").Value = "  -7.01%  "

$ws.Range("E29").Value = "  -2.42%  "

$ws.Range("E30").Value = "  -4.57%  "

$ws.Range("D31").Value = "172.82"
$ws.Range("E31").Value = "  -0.92%  "

$ws.Range("D32").Value = "36.84"
$ws.Range("E32").Value = "  -12.43%  "

$ws.Range("D33").Value = "20.84"
$ws.Range("E33").Value = "  -3.13%  "

$ws.Range("D34").Value = "0.0872"
$ws.Range("E34").Value = "  -4.94%  "

$ws.Range("D35").Value = "5.57"
$ws.Range("E35").Value = "  -2.90%  "

$ws.Range("D36").Value = "4.99"
$ws.Range("E36").Value = "  +8.59%  "

$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("D38").Value = "4.16"
$ws.Range("E38").Value = "  -2.51%  "

$ws.Range("D39").Value = "0.0368"
$ws.Range("E39").Value = "  -3.97%  "

$ws.Range("E40").Value = "  -3.77%  "

$ws.Range("D41").Value = "74.47"
$ws.Range("E41").Value = "  +2.53%  "

$ws.Range("D42").Value = "2.38"
$ws.Range("E42").Value = "  -7.20%  "

$ws.Range("D43").Value = "0.229"
$ws.Range("E43").Value = "  -4.79%  "

$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").Value = "12.29"
$ws.Range("E45").Value = "  -11.71%  "

$ws.Range("D46").Value = "1.29"
$ws.Range("E46").Value = "  -6.57%  "

$ws.Range("D47").Value = "5.37"
$ws.Range("E47").Value = "  -6.64%  "

$ws.Range("D48").Value = "1.74"
$ws.Range("E48").Value = "  +12.11%  "

$ws.Range("D49").Value = "1.27"
$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("E50").Value = "  -1.44%  "

$ws.Range("D51").Value = "101.28"
$ws.Range("E51").Value = "  -1.85%  "

# Restore the default cell style so we don't leave a stray
# number-format override on cells that previously had none.
$ws.Range("D2:E51").Style = "Normal"
